$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.336603
$ws.Cells.Item(2, 8).Value = 10.009809
$ws.Cells.Item(2, 9).Value = 0.2530037693731349
$ws.Cells.Item(2, 10).Value = 0.253003769373135
$ws.Cells.Item(2, 13).Value = 86.89540866666668
$ws.Cells.Item(2, 14).Value = 260.686226
$ws.Cells.Item(2, 15).Value = 0.319779657009892
$ws.Cells.Item(2, 16).Value = 0.3197796570098919
$ws.Cells.Item(2, 17).Value = 289.9354812434261
$ws.Cells.Item(2, 18).Value = 2609.419331190834
$ws.Cells.Item(2, 19).Value = 0.08090545859235089
$ws.Cells.Item(2, 20).Value = 0.08090545859235089
$ws.Cells.Item(3, 7).Value = 3.336603
$ws.Cells.Item(3, 8).Value = 10.009809
$ws.Cells.Item(3, 9).Value = 0.2530037693731349
$ws.Cells.Item(3, 10).Value = 0.253003769373135
$ws.Cells.Item(3, 15).Value = 0.1999969065479545
$ws.Cells.Item(3, 16).Value = 0.1999969065479545
$ws.Cells.Item(3, 17).Value = 181.331732885635
$ws.Cells.Item(3, 18).Value = 1631.985595970715
$ws.Cells.Item(3, 19).Value = 0.0505999712195991
$ws.Cells.Item(3, 20).Value = 0.05059997121959911
$ws.Cells.Item(4, 7).Value = 3.336603
$ws.Cells.Item(4, 8).Value = 10.009809
$ws.Cells.Item(4, 9).Value = 0.2530037693731349
$ws.Cells.Item(4, 10).Value = 0.253003769373135
$ws.Cells.Item(4, 13).Value = 60.92601633333334
$ws.Cells.Item(4, 14).Value = 182.778049
$ws.Cells.Item(4, 15).Value = 0.224210932487692
$ws.Cells.Item(4, 16).Value = 0.224210932487692
$ws.Cells.Item(4, 17).Value = 203.285928875849
$ws.Cells.Item(4, 18).Value = 1829.573359882641
$ws.Cells.Item(4, 19).Value = 0.05672621105405155
$ws.Cells.Item(4, 20).Value = 0.05672621105405155
$ws.Cells.Item(5, 7).Value = 3.336603
$ws.Cells.Item(5, 8).Value = 10.009809
$ws.Cells.Item(5, 9).Value = 0.2530037693731349
$ws.Cells.Item(5, 10).Value = 0.253003769373135
$ws.Cells.Item(5, 13).Value = 7.809668333333332
$ws.Cells.Item(5, 14).Value = 23.429005
$ws.Cells.Item(5, 15).Value = 0.02873998867505581
$ws.Cells.Item(5, 16).Value = 0.02873998867505581
$ws.Cells.Item(5, 17).Value = 26.057762790005
$ws.Cells.Item(5, 18).Value = 234.519865110045
$ws.Cells.Item(5, 19).Value = 0.00727132546653033
$ws.Cells.Item(5, 20).Value = 0.007271325466530332
$ws.Cells.Item(6, 7).Value = 3.336603
$ws.Cells.Item(6, 8).Value = 10.009809
$ws.Cells.Item(6, 9).Value = 0.2530037693731349
$ws.Cells.Item(6, 10).Value = 0.253003769373135
$ws.Cells.Item(6, 13).Value = 61.75795633333333
$ws.Cells.Item(6, 14).Value = 185.273869
$ws.Cells.Item(6, 15).Value = 0.2272725152794058
$ws.Cells.Item(6, 16).Value = 0.2272725152794058
$ws.Cells.Item(6, 17).Value = 206.061782375669
$ws.Cells.Item(6, 18).Value = 1854.556041381021
$ws.Cells.Item(6, 19).Value = 0.05750080304060307
$ws.Cells.Item(6, 20).Value = 0.05750080304060308
$ws.Cells.Item(7, 9).Value = 0.4389650376240971
$ws.Cells.Item(7, 10).Value = 0.4389650376240971
$ws.Cells.Item(7, 13).Value = 86.89540866666668
$ws.Cells.Item(7, 14).Value = 260.686226
$ws.Cells.Item(7, 15).Value = 0.319779657009892
$ws.Cells.Item(7, 16).Value = 0.3197796570098919
$ws.Cells.Item(7, 17).Value = 503.0420682977203
$ws.Cells.Item(7, 18).Value = 4527.378614679482
$ws.Cells.Item(7, 19).Value = 0.1403720891707681
$ws.Cells.Item(7, 20).Value = 0.140372089170768
$ws.Cells.Item(8, 9).Value = 0.4389650376240971
$ws.Cells.Item(8, 10).Value = 0.4389650376240971
$ws.Cells.Item(8, 15).Value = 0.1999969065479545
$ws.Cells.Item(8, 16).Value = 0.1999969065479545
$ws.Cells.Item(8, 19).Value = 0.08779164960752589
$ws.Cells.Item(8, 20).Value = 0.08779164960752588
$ws.Cells.Item(9, 9).Value = 0.4389650376240971
$ws.Cells.Item(9, 10).Value = 0.4389650376240971
$ws.Cells.Item(9, 13).Value = 60.92601633333334
$ws.Cells.Item(9, 14).Value = 182.778049
$ws.Cells.Item(9, 15).Value = 0.224210932487692
$ws.Cells.Item(9, 16).Value = 0.224210932487692
$ws.Cells.Item(9, 17).Value = 352.7038970151881
$ws.Cells.Item(9, 18).Value = 3174.335073136693
$ws.Cells.Item(9, 19).Value = 0.0984207604151936
$ws.Cells.Item(9, 20).Value = 0.0984207604151936
$ws.Cells.Item(10, 9).Value = 0.4389650376240971
$ws.Cells.Item(10, 10).Value = 0.4389650376240971
$ws.Cells.Item(10, 13).Value = 7.809668333333332
$ws.Cells.Item(10, 14).Value = 23.429005
$ws.Cells.Item(10, 15).Value = 0.02873998867505581
$ws.Cells.Item(10, 16).Value = 0.02873998867505581
$ws.Cells.Item(10, 17).Value = 45.21057868764277
$ws.Cells.Item(10, 18).Value = 406.8952081887849
$ws.Cells.Item(10, 19).Value = 0.012615850210062
$ws.Cells.Item(10, 20).Value = 0.012615850210062
$ws.Cells.Item(11, 9).Value = 0.4389650376240971
$ws.Cells.Item(11, 10).Value = 0.4389650376240971
$ws.Cells.Item(11, 13).Value = 61.75795633333333
$ws.Cells.Item(11, 14).Value = 185.273869
$ws.Cells.Item(11, 15).Value = 0.2272725152794058
$ws.Cells.Item(11, 16).Value = 0.2272725152794058
$ws.Cells.Item(11, 17).Value = 357.5200412133814
$ws.Cells.Item(11, 18).Value = 3217.680370920433
$ws.Cells.Item(11, 19).Value = 0.09976468822054756
$ws.Cells.Item(11, 20).Value = 0.09976468822054754
$ws.Cells.Item(12, 7).Value = 1.029432
$ws.Cells.Item(12, 8).Value = 3.088296
$ws.Cells.Item(12, 9).Value = 0.07805848532574147
$ws.Cells.Item(12, 10).Value = 0.07805848532574149
$ws.Cells.Item(12, 13).Value = 86.89540866666668
$ws.Cells.Item(12, 14).Value = 260.686226
$ws.Cells.Item(12, 15).Value = 0.319779657009892
$ws.Cells.Item(12, 16).Value = 0.3197796570098919
$ws.Cells.Item(12, 17).Value = 89.452914334544
$ws.Cells.Item(12, 18).Value = 805.076229010896
$ws.Cells.Item(12, 19).Value = 0.02496151566417729
$ws.Cells.Item(12, 20).Value = 0.02496151566417729
$ws.Cells.Item(13, 7).Value = 1.029432
$ws.Cells.Item(13, 8).Value = 3.088296
$ws.Cells.Item(13, 9).Value = 0.07805848532574147
$ws.Cells.Item(13, 10).Value = 0.07805848532574149
$ws.Cells.Item(13, 15).Value = 0.1999969065479545
$ws.Cells.Item(13, 16).Value = 0.1999969065479545
$ws.Cells.Item(13, 17).Value = 55.94572936844
$ws.Cells.Item(13, 18).Value = 503.51156431596
$ws.Cells.Item(13, 19).Value = 0.0156114555949672
$ws.Cells.Item(13, 20).Value = 0.0156114555949672
$ws.Cells.Item(14, 7).Value = 1.029432
$ws.Cells.Item(14, 8).Value = 3.088296
$ws.Cells.Item(14, 9).Value = 0.07805848532574147
$ws.Cells.Item(14, 10).Value = 0.07805848532574149
$ws.Cells.Item(14, 13).Value = 60.92601633333334
$ws.Cells.Item(14, 14).Value = 182.778049
$ws.Cells.Item(14, 15).Value = 0.224210932487692
$ws.Cells.Item(14, 16).Value = 0.224210932487692
$ws.Cells.Item(14, 17).Value = 62.71919084605599
$ws.Cells.Item(14, 18).Value = 564.4727176145039
$ws.Cells.Item(14, 19).Value = 0.01750156578346132
$ws.Cells.Item(14, 20).Value = 0.01750156578346132
$ws.Cells.Item(15, 7).Value = 1.029432
$ws.Cells.Item(15, 8).Value = 3.088296
$ws.Cells.Item(15, 9).Value = 0.07805848532574147
$ws.Cells.Item(15, 10).Value = 0.07805848532574149
$ws.Cells.Item(15, 13).Value = 7.809668333333332
$ws.Cells.Item(15, 14).Value = 23.429005
$ws.Cells.Item(15, 15).Value = 0.02873998867505581
$ws.Cells.Item(15, 16).Value = 0.02873998867505581
$ws.Cells.Item(15, 17).Value = 8.039522491719998
$ws.Cells.Item(15, 18).Value = 72.35570242547998
$ws.Cells.Item(15, 19).Value = 0.00224339998425382
$ws.Cells.Item(15, 20).Value = 0.002243399984253821
$ws.Cells.Item(16, 7).Value = 1.029432
$ws.Cells.Item(16, 8).Value = 3.088296
$ws.Cells.Item(16, 9).Value = 0.07805848532574147
$ws.Cells.Item(16, 10).Value = 0.07805848532574149
$ws.Cells.Item(16, 13).Value = 61.75795633333333
$ws.Cells.Item(16, 14).Value = 185.273869
$ws.Cells.Item(16, 15).Value = 0.2272725152794058
$ws.Cells.Item(16, 16).Value = 0.2272725152794058
$ws.Cells.Item(16, 17).Value = 63.57561650413599
$ws.Cells.Item(16, 18).Value = 572.1805485372239
$ws.Cells.Item(16, 19).Value = 0.01774054829888186
$ws.Cells.Item(16, 20).Value = 0.01774054829888186
$ws.Cells.Item(17, 7).Value = 1.675087666666667
$ws.Cells.Item(17, 8).Value = 5.025263
$ws.Cells.Item(17, 9).Value = 0.1270164576658104
$ws.Cells.Item(17, 10).Value = 0.1270164576658104
$ws.Cells.Item(17, 13).Value = 86.89540866666668
$ws.Cells.Item(17, 14).Value = 260.686226
$ws.Cells.Item(17, 15).Value = 0.319779657009892
$ws.Cells.Item(17, 16).Value = 0.3197796570098919
$ws.Cells.Item(17, 17).Value = 145.5574273474931
$ws.Cells.Item(17, 18).Value = 1310.016846127438
$ws.Cells.Item(17, 19).Value = 0.04061727926698432
$ws.Cells.Item(17, 20).Value = 0.04061727926698431
$ws.Cells.Item(18, 7).Value = 1.675087666666667
$ws.Cells.Item(18, 8).Value = 5.025263
$ws.Cells.Item(18, 9).Value = 0.1270164576658104
$ws.Cells.Item(18, 10).Value = 0.1270164576658104
$ws.Cells.Item(18, 15).Value = 0.1999969065479545
$ws.Cells.Item(18, 16).Value = 0.1999969065479545
$ws.Cells.Item(18, 17).Value = 91.03466889288944
$ws.Cells.Item(18, 18).Value = 819.3120200360049
$ws.Cells.Item(18, 19).Value = 0.02540289861384131
$ws.Cells.Item(18, 20).Value = 0.0254028986138413
$ws.Cells.Item(19, 7).Value = 1.675087666666667
$ws.Cells.Item(19, 8).Value = 5.025263
$ws.Cells.Item(19, 9).Value = 0.1270164576658104
$ws.Cells.Item(19, 10).Value = 0.1270164576658104
$ws.Cells.Item(19, 13).Value = 60.92601633333334
$ws.Cells.Item(19, 14).Value = 182.778049
$ws.Cells.Item(19, 15).Value = 0.224210932487692
$ws.Cells.Item(19, 16).Value = 0.224210932487692
$ws.Cells.Item(19, 17).Value = 102.0564185390986
$ws.Cells.Item(19, 18).Value = 918.507766851887
$ws.Cells.Item(19, 19).Value = 0.02847847841453481
$ws.Cells.Item(19, 20).Value = 0.0284784784145348
$ws.Cells.Item(20, 7).Value = 1.675087666666667
$ws.Cells.Item(20, 8).Value = 5.025263
$ws.Cells.Item(20, 9).Value = 0.1270164576658104
$ws.Cells.Item(20, 10).Value = 0.1270164576658104
$ws.Cells.Item(20, 13).Value = 7.809668333333332
$ws.Cells.Item(20, 14).Value = 23.429005
$ws.Cells.Item(20, 15).Value = 0.02873998867505581
$ws.Cells.Item(20, 16).Value = 0.02873998867505581
$ws.Cells.Item(20, 17).Value = 13.08187910592389
$ws.Cells.Item(20, 18).Value = 117.736911953315
$ws.Cells.Item(20, 19).Value = 0.003650451554861097
$ws.Cells.Item(20, 20).Value = 0.003650451554861097
$ws.Cells.Item(21, 7).Value = 1.675087666666667
$ws.Cells.Item(21, 8).Value = 5.025263
$ws.Cells.Item(21, 9).Value = 0.1270164576658104
$ws.Cells.Item(21, 10).Value = 0.1270164576658104
$ws.Cells.Item(21, 13).Value = 61.75795633333333
$ws.Cells.Item(21, 14).Value = 185.273869
$ws.Cells.Item(21, 15).Value = 0.2272725152794058
$ws.Cells.Item(21, 16).Value = 0.2272725152794058
$ws.Cells.Item(21, 17).Value = 103.4499909725052
$ws.Cells.Item(21, 18).Value = 931.0499187525469
$ws.Cells.Item(21, 19).Value = 0.0288673498155889
$ws.Cells.Item(21, 20).Value = 0.0288673498155889
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 1.357782666666667
$ws.Cells.Item(22, 8).Value = 4.073348
$ws.Cells.Item(22, 9).Value = 0.102956250011216
$ws.Cells.Item(22, 10).Value = 0.102956250011216
$ws.Cells.Item(22, 13).Value = 86.89540866666668
$ws.Cells.Item(22, 14).Value = 260.686226
$ws.Cells.Item(22, 15).Value = 0.319779657009892
$ws.Cells.Item(22, 16).Value = 0.3197796570098919
$ws.Cells.Item(22, 17).Value = 117.9850797005165
$ws.Cells.Item(22, 18).Value = 1061.865717304648
$ws.Cells.Item(22, 19).Value = 0.03292331431561135
$ws.Cells.Item(22, 20).Value = 0.03292331431561134
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 1.357782666666667
$ws.Cells.Item(23, 8).Value = 4.073348
$ws.Cells.Item(23, 9).Value = 0.102956250011216
$ws.Cells.Item(23, 10).Value = 0.102956250011216
$ws.Cells.Item(23, 15).Value = 0.1999969065479545
$ws.Cells.Item(23, 16).Value = 0.1999969065479545
$ws.Cells.Item(23, 17).Value = 73.79034419999779
$ws.Cells.Item(23, 18).Value = 664.11309779998
$ws.Cells.Item(23, 19).Value = 0.02059093151202102
$ws.Cells.Item(23, 20).Value = 0.02059093151202102
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 1.357782666666667
$ws.Cells.Item(24, 8).Value = 4.073348
$ws.Cells.Item(24, 9).Value = 0.102956250011216
$ws.Cells.Item(24, 10).Value = 0.102956250011216
$ws.Cells.Item(24, 13).Value = 60.92601633333334
$ws.Cells.Item(24, 14).Value = 182.778049
$ws.Cells.Item(24, 15).Value = 0.224210932487692
$ws.Cells.Item(24, 16).Value = 0.224210932487692
$ws.Cells.Item(24, 17).Value = 82.72428892645023
$ws.Cells.Item(24, 18).Value = 744.5186003380521
$ws.Cells.Item(24, 19).Value = 0.0230839168204507
$ws.Cells.Item(24, 20).Value = 0.0230839168204507
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 1.357782666666667
$ws.Cells.Item(25, 8).Value = 4.073348
$ws.Cells.Item(25, 9).Value = 0.102956250011216
$ws.Cells.Item(25, 10).Value = 0.102956250011216
$ws.Cells.Item(25, 13).Value = 7.809668333333332
$ws.Cells.Item(25, 14).Value = 23.429005
$ws.Cells.Item(25, 15).Value = 0.02873998867505581
$ws.Cells.Item(25, 16).Value = 0.02873998867505581
$ws.Cells.Item(25, 17).Value = 10.60383229541555
$ws.Cells.Item(25, 18).Value = 95.43449065873999
$ws.Cells.Item(25, 19).Value = 0.002958961459348564
$ws.Cells.Item(25, 20).Value = 0.002958961459348564
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 1.357782666666667
$ws.Cells.Item(26, 8).Value = 4.073348
$ws.Cells.Item(26, 9).Value = 0.102956250011216
$ws.Cells.Item(26, 10).Value = 0.102956250011216
$ws.Cells.Item(26, 13).Value = 61.75795633333333
$ws.Cells.Item(26, 14).Value = 185.273869
$ws.Cells.Item(26, 15).Value = 0.2272725152794058
$ws.Cells.Item(26, 16).Value = 0.2272725152794058
$ws.Cells.Item(26, 17).Value = 83.85388263815689
$ws.Cells.Item(26, 18).Value = 754.684943743412
$ws.Cells.Item(26, 19).Value = 0.02339912590378443
$ws.Cells.Item(26, 20).Value = 0.02339912590378442
